$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix photo paths: .jpg -> .webp in the "Photo" column ---
$used = $ws.UsedRange
$lastRow = $used.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value()
    if ($val -ne $null -and $val -like "*.jpg") {
        $cell.Value = $val -replace "\.jpg$", ".webp"
    }
}

# --- 2) Header font color: theme color -> explicit black RGB ---
$header = $ws.Range("A1:J1")
$header.Font.Color = 0

# --- 3) Column width adjustments for Designation (B) and Photo (C) ---
$ws.Columns("B").ColumnWidth = 20.576428571428572
$ws.Columns("C").ColumnWidth = 50.14785714285715
